$d = $word.ActiveDocument

# Locate the whole sentence first, to scope our inner search for "and" to this specific
# paragraph (the word "and" occurs many other times throughout the document).
$sentence = $d.Content
$found = $sentence.Find.Execute(
    "Prepare a query that matches the actor’s request and retrieve the data from",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $found) {
    throw "Could not find target sentence"
}
$sentenceStart = $sentence.Start
$sentenceEnd = $sentence.End

# Now search for "and" within that scoped range only, so we target exactly the right word.
$scoped = $d.Range($sentenceStart, $sentenceEnd)
$found2 = $scoped.Find.Execute("and", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $found2) {
    throw "Could not find 'and' within sentence"
}
$andStart = $scoped.Start
$andEnd = $scoped.End

# Replace "and" with "or" in place.
$d.Range($andStart, $andEnd).Delete()
$insertionPoint = $d.Range($andStart, $andStart)
$insertionPoint.InsertAfter("or")

# The insertion merged "or" into the surrounding run (identical formatting). Force Word to
# split it into its own run -- matching the target -- by toggling a character property on
# just the new "or" text and then reverting it.
$orRange = $d.Range($andStart, $andStart + 2)
$orRange.Bold = 1
$orRange2 = $d.Range($andStart, $andStart + 2)
$orRange2.Bold = 0
